$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.810.75'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '1.642.89'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('E4').Value = '  -0.40%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.94'
$ws.Range('E5').Value = '  -0.38%  '
$ws.Range('E6').Value = '  +0.44%  '
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('E8').Value = '  +1.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0619'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.67'
$ws.Range('E10').Value = '  +3.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0844'
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('D12').Value = '1.874.19'
$ws.Range('E12').Value = '  +0.39%  '
$ws.Range('D13').Value = '1.653.76'
$ws.Range('E13').Value = '  +1.50%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.13'
$ws.Range('E14').Value = '  -0.06%  '
$ws.Range('E15').Value = '  +0.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.11'
$ws.Range('E16').Value = '  +2.58%  '
$ws.Range('D17').Value = '26.846.64'
$ws.Range('E17').Value = '  +0.59%  '
$ws.Range('D18').Value = '0.0₃0729'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '217.63'
$ws.Range('E19').Value = '  +3.06%  '
$ws.Range('E20').Value = '  -0.44%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.39'
$ws.Range('E21').Value = '  +1.31%  '
$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.63'
$ws.Range('E22').Value = '  +7.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.45'
$ws.Range('E23').Value = '  +5.87%  '
$ws.Range('E24').Value = '  -0.90%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.94'
$ws.Range('E25').Value = '  -0.49%  '
$ws.Range('E26').Value = '  -0.52%  '
$ws.Range('E27').Value = '  +3.64%  '
$ws.Range('E28').Value = '  +0.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.82'
$ws.Range('E29').Value = '  +1.75%  '
$ws.Range('E30').Value = '  +1.62%  '
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.37'
$ws.Range('E32').Value = '  +0.78%  '
$ws.Range('E33').Value = '  +0.32%  '
$ws.Range('E34').Value = '  +1.59%  '
$ws.Range('D36').Value = '1.245.22'
$ws.Range('E36').Value = '  -1.86%  '
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.534'
$ws.Range('E38').Value = '  +1.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.831'
$ws.Range('E39').Value = '  +3.35%  '
$ws.Range('E40').Value = '  -0.42%  '
$ws.Range('E41').Value = '  +0.50%  '
$ws.Range('E42').Value = '  +1.52%  '
$ws.Range('D43').Value = '1.785.93'
$ws.Range('E43').Value = '  +0.57%  '
$ws.Range('E44').Value = '  -3.95%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '60.85'
$ws.Range('E45').Value = '  +1.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.55'
$ws.Range('E46').Value = '  +0.33%  '
$ws.Range('E47').Value = '  +0.80%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0514'
$ws.Range('E48').Value = '  -0.84%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0971'
$ws.Range('E49').Value = '  +1.23%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.51'
$ws.Range('E50').Value = '  -0.13%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.406'
$ws.Range('E51').Value = '  -0.26%  '
